$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The app now outputs a summary sheet of uploaded items: the per-row
# "sub component" (column C) is dropped, and column B now carries a
# per-item comment/status instead of the old self-subcomponent test text.

# Row 7 header: drop the "child" sub-component column heading.
$ws.Range("C7").ClearContents()

# Row 8: same serial number, new comment, no sub-component.
$ws.Range("B8").Value = "testing output sheets"
$ws.Range("C8").ClearContents()

# Row 9: new item uploaded, flagged "new", no sub-component.
$ws.Range("A9").Value = "SN990102"
$ws.Range("B9").Value = "new"
$ws.Range("C9").ClearContents()

# Row 10: another new item uploaded, flagged "new", no sub-component.
$ws.Range("A10").Value = "SN990103"
$ws.Range("B10").Value = "new"
$ws.Range("C10").ClearContents()

# Selection moved from D9 to C9.
$ws.Range("C9").Select()
